$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels: columns F/G are brand-new (height, weight); the old
# "fantasy points" header slides from E1 over to G1.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# Give the two new header cells the same look as the rest of row 1
# (bold / centered / bordered header style).
$ws.Range("D1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Per-row data: move the old fantasy-points value from column E into the
# new column G, then fill E (height) and F (weight) with the scraped
# constants.
$fantasyPoints = @(0, 1.9, 4.4, 1.2, 6.5, 0, 1.3, 0, 3.8, 1.2, 8, 0, 1.5, 0.8)

for ($i = 0; $i -lt $fantasyPoints.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $fantasyPoints[$i]
    $ws.Cells.Item($row, 5).Value = 6.5
    $ws.Cells.Item($row, 6).Value = 255
}

Write-Output "done"
